$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.141.94"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "3.462.19"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "576.79"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "160.92"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.23%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.612"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +6.63%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "3.465.02"
$ws.Range("E9").Value = "  +1.77%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.26"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  +2.01%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.454"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").Value = "4.057.47"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("E14").Value = "  +0.74%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000193"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.78%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "28.31"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "65.099.28"
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "3.466.62"
$ws.Range("E18").Value = "  +1.89%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.49"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.37"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.55%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "382.79"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "8.18"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.30%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.556"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.85%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "73.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.21%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  +1.50%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.07"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.18%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.178"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.18%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +10.14%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.18"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.65%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.05"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.60%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.70"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.55%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.30"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.29%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.62"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +11.87%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "161.56"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("E37").Value = "  +5.35%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0782"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.85%  "
$ws.Range("D39").Value = "2.917.21"
$ws.Range("E39").Value = "  +0.74%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.84"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.89%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.95"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "4.70"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +7.49%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "43.19"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0320"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.29%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.780"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.50%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "26.02"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +12.23%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "325.19"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +11.84%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.10"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("E49").Value = "  +3.10%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.878"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.70%  "
$ws.Range("E51").Value = "  +0.28%  "
